$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 84
$ws.Range("A84").Value = "Réalisation"
$ws.Range("A84").WrapText = $true
$ws.Range("B84").Value = "Je continue d'extraire les informations du fichier GPX, je recontre quelques difficultés car le module nodejs que j'utilise (gpx-parse) ne possède pas énormément de documentation "
$ws.Range("C84").Value = 0.75
$ws.Range("D84").Value = 43550
$ws.Rows.Item(84).RowHeight = 105

# Row 85
$ws.Range("A85").Value = "Réalisation"
$ws.Range("A85").WrapText = $true
$ws.Range("B85").Value = "Continuation de l'extraction des données du GPX, j'arrive désormais a créer des enregistrements dans la table des positions. J'ai commencé à ajouté les champs qui vont être calculés à la base de données"
$ws.Range("C85").Value = 1.5
$ws.Range("D85").Value = 43551
$ws.Rows.Item(85).RowHeight = 105

# Row 86
$ws.Range("A86").Value = "Réalisation"
$ws.Range("A86").WrapText = $true
$ws.Range("B86").Value = "J'effectue réalisé les calculs pour les nouveaux champs de l'entité d'activité. "
$ws.Range("C86").Value = 3
$ws.Range("D86").Value = 43551
$ws.Range("E86").Value = "J'ai rencontré des problèmes avec le fichier GPX fournit par M. Glassey, un problème lors de l'exportation de celui-ci à du survenir, car les timestamps de celui-ci n'était pas correct, ce qui faussait mes résultats."
$ws.Range("E86").WrapText = $true
$ws.Rows.Item(86).RowHeight = 60

# Row 87
$ws.Range("A87").Value = "Réalisation"
$ws.Range("A87").WrapText = $true
$ws.Range("B87").Value = "Correction d'un bug lors de la création d'activtié qui empêchait la création sans fichier gpx"
$ws.Range("C87").Value = 1
$ws.Range("D87").Value = 43552
$ws.Rows.Item(87).RowHeight = 60

# Row 88
$ws.Range("A88").Value = "Gestion de projet"
$ws.Range("A88").WrapText = $true
$ws.Range("B88").Value = "Mise à jour de Trello, je prends contact avec M. Glassey pour planifier la tâche suivante. Je souhaiterai travailler sur l'authentification par token"
$ws.Range("C88").Value = 0.5
$ws.Range("D88").Value = 43552
$ws.Rows.Item(88).RowHeight = 75

# Row 89
$ws.Range("A89").Value = "Analyse"
$ws.Range("A89").WrapText = $true
$ws.Range("B89").Value = "Je me rends compte d'un besoin d'implémenter une gestion de token non-JWT pour l'interface web"
$ws.Range("C89").Value = 0.5
$ws.Range("D89").Value = 43552
$ws.Rows.Item(89).RowHeight = 60

# Row 90
$ws.Range("A90").Value = "Conception"
$ws.Range("A90").WrapText = $true
$ws.Range("B90").Value = "J'image la forme que peut prendre ce nouveau composant au sein de l'API. Je pense ajouté une nouvelle table gérant les token de session"
$ws.Range("C90").Value = 1
$ws.Range("D90").Value = 43552
$ws.Rows.Item(90).RowHeight = 75

# Row 91
$ws.Range("A91").Value = "Réalisation"
$ws.Range("A91").WrapText = $true
$ws.Range("B91").Value = "Mise à jour du MLD, mise à jour du schéma de la base de données"
$ws.Range("C91").Value = 1
$ws.Range("D91").Value = 43552
$ws.Rows.Item(91).RowHeight = 45

# Row 92
$ws.Range("A92").Value = "Gestion de projet"
$ws.Range("A92").WrapText = $true
$ws.Range("B92").Value = "Mise à jour du Journal de bord"
$ws.Range("C92").Value = 0.5
$ws.Range("D92").Value = 76424

$ws.Range("B86").Select()
